$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B:B").HorizontalAlignment = -4131
$ws.Columns("C:C").Insert()
$ws.Range("B1").Value = "nmff_mrn"
$ws.Range("C1").Value = "nmh_mrn"
